$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Admin's password / hash changed ("chore: changed default password")
$ws.Range("C2").Value = '$2b$12$QxPEm5z1OAJzc0SFHvqKNecF7nq4DkeciWCVAPcXZgkTE3NAlU8JS'
$ws.Range("D2").Value = '6a1bbc448a668c31e0cfe425b018f299b0d9dcb666ebb5ec1819e276d5ec9550'

# Highlight fulano's hash cell with an underline
$ws.Range("D3").Font.Underline = $true

# Move the active selection
[void]$ws.Range("C7").Select()
